$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: direct assignment is safe.
$ws.Range("D2").Value = '29.320.22'
$ws.Range("E2").Value = '  -0.52%  '
$ws.Range("D3").Value = '1.845.62'
$ws.Range("E3").Value = '  -0.37%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("E6").Value = '  -0.34%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("E9").Value = '  -0.92%  '
$ws.Range("E10").Value = '  -0.60%  '
$ws.Range("E11").Value = '  -0.06%  '
$ws.Range("D12").Value = '1.847.05'
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("E14").Value = '  -0.52%  '
$ws.Range("E15").Value = '  -2.06%  '
$ws.Range("E16").Value = '  -0.90%  '
$ws.Range("E17").Value = '  -1.30%  '
$ws.Range("D18").Value = '29.357.29'
$ws.Range("E18").Value = '  -0.51%  '
$ws.Range("E19").Value = '  -0.33%  '
$ws.Range("E20").Value = '  -1.15%  '
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("E22").Value = '  -0.34%  '
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("E24").Value = '  +1.25%  '
$ws.Range("E25").Value = '  +0.91%  '
$ws.Range("E26").Value = '  +0.41%  '
$ws.Range("E27").Value = '  -0.31%  '
$ws.Range("E28").Value = '  +6.91%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("E30").Value = '  -0.17%  '
$ws.Range("E31").Value = '  -0.76%  '
$ws.Range("E32").Value = '  +0.53%  '
$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("E33").Value = '  -1.71%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("E34").Value = '  -1.19%  '
$ws.Range("E35").Value = '  -1.69%  '
$ws.Range("E36").Value = '  -0.34%  '
$ws.Range("E37").Value = '  +2.06%  '
$ws.Range("D38").Value = '1.239.96'
$ws.Range("E38").Value = '  +1.45%  '
$ws.Range("E39").Value = '  -2.28%  '
$ws.Range("E40").Value = '  -2.83%  '
$ws.Range("E41").Value = '  -0.87%  '
$ws.Range("E42").Value = '  -0.27%  '
$ws.Range("D43").Value = '2.006.19'
$ws.Range("E43").Value = '  -0.91%  '
$ws.Range("E44").Value = '  -0.30%  '
$ws.Range("E45").Value = '  -1.51%  '
$ws.Range("E46").Value = '  -0.50%  '
$ws.Range("B47").Value = 'TheSandbox'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("E47").Value = '  -0.80%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("E48").Value = '  -3.40%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("E49").Value = '  -0.56%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("E50").Value = '  +0.35%  '
$ws.Range("E51").Value = '  -0.64%  '

# Numeric-looking strings (e.g. "0.9994") must be forced to stay as text,
# matching the source workbook where these Price cells are inline strings,
# not numbers. We briefly mark the cell as Text, assign, then restore the
# default "Normal" style so no stray formatting is left behind.
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '0.9994'
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '240.27'
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.6285'
$cell.Style = "Normal"
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.9995'
$cell.Style = "Normal"
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.07576'
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.2911'
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '24.56'
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.07746'
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '5.009'
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '0.6779'
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.00001048'
$cell.Style = "Normal"
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '83.02'
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '6.102'
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '229.07'
$cell.Style = "Normal"
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '0.9997'
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '7.439'
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '158.92'
$cell.Style = "Normal"
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '0.1396'
$cell.Style = "Normal"
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '8.441'
$cell.Style = "Normal"
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '1.418'
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '1.470'
$cell.Style = "Normal"
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '0.05674'
$cell.Style = "Normal"
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '4.107'
$cell.Style = "Normal"
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '4.066'
$cell.Style = "Normal"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '1.821'
$cell.Style = "Normal"
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.153'
$cell.Style = "Normal"
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '0.6961'
$cell.Style = "Normal"
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '2.580'
$cell.Style = "Normal"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '0.01833'
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '2.719'
$cell.Style = "Normal"
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '6.371'
$cell.Style = "Normal"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.9013'
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.9988'
$cell.Style = "Normal"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '101.33'
$cell.Style = "Normal"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '65.49'
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '7.109'
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.3992'
$cell.Style = "Normal"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.00000000116'
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '8.998'
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.1149'
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '1.676'
$cell.Style = "Normal"
